$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D holds price strings that look numeric (e.g. "61.152.88", "1.01").
# Force text storage (as in the source file, t="inlineStr") by setting the
# number format to Text before assigning the value - mirrors typing the value
# into a cell that has already been formatted as Text in Excel.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '61.152.88'
$ws.Range('E2').Value = '  +11.17%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.623.31'
$ws.Range('E3').Value = '  +12.38%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '517.67'
$ws.Range('E5').Value = '  +9.19%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '162.86'
$ws.Range('E6').Value = '  +12.40%  '
$ws.Range('E7').Value = '  +2.53%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.992'
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.684.47'
$ws.Range('E9').Value = '  +14.87%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.23'
$ws.Range('E10').Value = '  +14.70%  '
$ws.Range('E11').Value = '  +12.40%  '
$ws.Range('E12').Value = '  +8.58%  '
$ws.Range('E13').Value = '  +1.77%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.091.73'
$ws.Range('E14').Value = '  +12.65%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '61.115.46'
$ws.Range('E15').Value = '  +10.98%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '22.73'
$ws.Range('E16').Value = '  +14.53%  '
$ws.Range('E17').Value = '  +11.04%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.671.91'
$ws.Range('E18').Value = '  +14.11%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.88'
$ws.Range('E19').Value = '  +6.88%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '354.13'
$ws.Range('E20').Value = '  +13.19%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.73'
$ws.Range('E21').Value = '  +13.05%  '
$ws.Range('E22').Value = '  +11.13%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.996'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '60.90'
$ws.Range('E24').Value = '  +8.68%  '
$ws.Range('E25').Value = '  +10.01%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.171'
$ws.Range('E26').Value = '  +12.06%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.782.92'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.993'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('E29').Value = '  +19.34%  '
$ws.Range('E30').Value = '  +9.71%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.997'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('E32').Value = '  +9.55%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '158.45'
$ws.Range('E33').Value = '  +8.54%  '
$ws.Range('E34').Value = '  +9.71%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.74'
$ws.Range('E35').Value = '  +13.10%  '
$ws.Range('E36').Value = '  +13.90%  '
$ws.Range('E37').Value = '  +13.17%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.886'
$ws.Range('E38').Value = '  +10.32%  '
$ws.Range('E39').Value = '  +15.89%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.858'
$ws.Range('E40').Value = '  +39.19%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.83'
$ws.Range('E41').Value = '  +14.10%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '304.87'
$ws.Range('E42').Value = '  +23.30%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '36.21'
$ws.Range('E43').Value = '  +8.06%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.647'
$ws.Range('E44').Value = '  +12.71%  '
$ws.Range('E45').Value = '  +14.85%  '
$ws.Range('E46').Value = '  +1.49%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '20.43'
$ws.Range('E47').Value = '  +23.16%  '
$ws.Range('E48').Value = '  +16.00%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.987'
$ws.Range('E49').Value = '  -1.24%  '
$ws.Range('E51').Value = '  +20.52%  '
